$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend header row (H1:J1) with new columns, copying header style from G1 ---
$ws.Range("G1").Copy()
$ws.Range("H1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("H1").Value = "antecedent_len"
$ws.Range("I1").Value = "consequent_len"
$ws.Range("J1").Value = "combo_len"

# --- Row 2 ---
$ws.Range("A2").Value = "frozenset({'ACCESS_NETWORK_STATE'})"
$ws.Range("B2").Value = "frozenset({'INTERNET'})"
$ws.Range("C2").Value = 0.988095238095238
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = "inf"
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 2

# --- Row 3 ---
$ws.Range("A3").Value = "frozenset({'WAKE_LOCK'})"
$ws.Range("B3").Value = "frozenset({'ACCESS_NETWORK_STATE'})"
$ws.Range("C3").Value = 0.9285714285714286
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 1.012048192771084
$ws.Range("F3").Value = 0.01105442176870741
$ws.Range("G3").Value = "inf"
$ws.Range("H3").Value = 1
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 2

# --- Row 4 (new) ---
$ws.Range("A4").Value = "frozenset({'WAKE_LOCK', 'INTERNET'})"
$ws.Range("B4").Value = "frozenset({'ACCESS_NETWORK_STATE'})"
$ws.Range("C4").Value = 0.9285714285714286
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1.012048192771084
$ws.Range("F4").Value = 0.01105442176870741
$ws.Range("G4").Value = "inf"
$ws.Range("H4").Value = 2
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 3

# --- Row 5 (new) ---
$ws.Range("A5").Value = "frozenset({'ACCESS_NETWORK_STATE'})"
$ws.Range("B5").Value = "frozenset({'WAKE_LOCK'})"
$ws.Range("C5").Value = 0.9285714285714286
$ws.Range("D5").Value = 0.9397590361445785
$ws.Range("E5").Value = 1.012048192771084
$ws.Range("F5").Value = 0.01105442176870741
$ws.Range("G5").Value = 1.185714285714285
$ws.Range("H5").Value = 1
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 2

# --- Row 6 (new) ---
$ws.Range("A6").Value = "frozenset({'ACCESS_NETWORK_STATE'})"
$ws.Range("B6").Value = "frozenset({'WAKE_LOCK', 'INTERNET'})"
$ws.Range("C6").Value = 0.9285714285714286
$ws.Range("D6").Value = 0.9397590361445785
$ws.Range("E6").Value = 1.012048192771084
$ws.Range("F6").Value = 0.01105442176870741
$ws.Range("G6").Value = 1.185714285714285
$ws.Range("H6").Value = 1
$ws.Range("I6").Value = 2
$ws.Range("J6").Value = 3

# --- Row 7 (new) ---
$ws.Range("A7").Value = "frozenset({'INTERNET'})"
$ws.Range("B7").Value = "frozenset({'WAKE_LOCK'})"
$ws.Range("C7").Value = 0.9285714285714286
$ws.Range("D7").Value = 0.9285714285714286
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 1
$ws.Range("I7").Value = 1
$ws.Range("J7").Value = 2
